$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B45").Value = "Some Useful Resources"
